$d = $word.ActiveDocument

# Remove "Standard" from the document title
$d.Content.Find.Execute("Square One Standard Statement of Work", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Square One Statement of Work", 2)

# Remove "Standard" from the contractor terms reference
$d.Content.Find.Execute("Square One Standard Contractor Terms", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Square One Contractor Terms", 2)
